# Weekly update: prepend a new week's Coliflor records (Primera/Segunda)
# for Femacal de La Calera, shifting the existing data down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 1011-1012, pushing old rows 1011.. down to 1013..
$ws.Range("A1011:A1012").EntireRow.Insert()

# New row 1011 - "Primera" quality record for the new week (Fecha 45106)
$ws.Cells.Item(1011, 1).Value = 3
$ws.Cells.Item(1011, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(1011, 3).Value = "Coquimbo"
$ws.Cells.Item(1011, 4).Value = 45106
$ws.Cells.Item(1011, 5).Value = 5
$ws.Cells.Item(1011, 6).Value = 100112008
$ws.Cells.Item(1011, 7).Value = "Coliflor"
$ws.Cells.Item(1011, 8).Value = "Sin especificar"
$ws.Cells.Item(1011, 9).Value = "Primera"
$ws.Cells.Item(1011, 10).Value = 3500
$ws.Cells.Item(1011, 11).Value = 800
$ws.Cells.Item(1011, 12).Value = 850
$ws.Cells.Item(1011, 13).Value = 823
$ws.Cells.Item(1011, 14).Value = "$/unidad"
$ws.Cells.Item(1011, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(1011, 16).Value = 823
$ws.Cells.Item(1011, 17).Value = 1
$ws.Cells.Item(1011, 18).Value = "Hortaliza"

# New row 1012 - "Segunda" quality record for the new week (Fecha 45106)
$ws.Cells.Item(1012, 1).Value = 3
$ws.Cells.Item(1012, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(1012, 3).Value = "Coquimbo"
$ws.Cells.Item(1012, 4).Value = 45106
$ws.Cells.Item(1012, 5).Value = 5
$ws.Cells.Item(1012, 6).Value = 100112008
$ws.Cells.Item(1012, 7).Value = "Coliflor"
$ws.Cells.Item(1012, 8).Value = "Sin especificar"
$ws.Cells.Item(1012, 9).Value = "Segunda"
$ws.Cells.Item(1012, 10).Value = 1800
$ws.Cells.Item(1012, 11).Value = 700
$ws.Cells.Item(1012, 12).Value = 700
$ws.Cells.Item(1012, 13).Value = 700
$ws.Cells.Item(1012, 14).Value = "$/unidad"
$ws.Cells.Item(1012, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(1012, 16).Value = 700
$ws.Cells.Item(1012, 17).Value = 1
$ws.Cells.Item(1012, 18).Value = "Hortaliza"
